$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.461.72'
$ws.Range("E2").Value = '  -1.94%  '
$ws.Range("D3").Value = '3.173.83'
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.74%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.172.12'
$ws.Range("E8").Value = '  -3.44%  '
$ws.Range("E9").Value = '  -2.09%  '
$ws.Range("E10").Value = '  -5.18%  '
$ws.Range("E11").Value = '  -2.56%  '
$ws.Range("E12").Value = '  -3.22%  '
$ws.Range("E13").Value = '  -4.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.63%  '
$ws.Range("E15").Value = '  -3.58%  '
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("D17").Value = '3.173.22'
$ws.Range("E17").Value = '  -3.64%  '
$ws.Range("D18").Value = '62.456.69'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '454.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("E22").Value = '  -3.63%  '
$ws.Range("E23").Value = '  -4.75%  '
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.53%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -3.12%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  -6.13%  '
$ws.Range("E30").Value = '  -4.70%  '
$ws.Range("E31").Value = '  -7.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.53%  '
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("E34").Value = '  -5.99%  '
$ws.Range("E35").Value = '  -5.95%  '
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.61%  '
$ws.Range("D38").Value = '0.0₃0698'
$ws.Range("E38").Value = '  -5.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0383'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '395.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.87%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.34%  '
$ws.Range("D44").Value = '2.795.05'
$ws.Range("E44").Value = '  -8.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.250'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.50%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.110'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.83%  '
